$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet: update "想去人数" (interest count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 36
$wsExhibit.Range("F5").Value = 4933
$wsExhibit.Range("F6").Value = 168
$wsExhibit.Range("F7").Value = 74
$wsExhibit.Range("F8").Value = 272

# "全部类型" (All types) sheet: same rows duplicated here, keep them in sync
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 36
$wsAll.Range("F9").Value = 4933
$wsAll.Range("F10").Value = 168
$wsAll.Range("F11").Value = 74
$wsAll.Range("F13").Value = 272
